# NUvention Web Scenario Assignment - wording fixes.
# Splits specific runs into multiple runs with corrected wording, matching
# the target OOXML diff. Each edit:
#   1. Locates the old (full) run text inside the shape's TextRange.
#   2. Overwrites that whole span with the concatenation of the new parts
#      (so the visible/plain text is correct even before the split).
#   3. Re-assigns each individual part onto its own character sub-range so
#      the engine materialises separate <a:r> runs (cloning the original
#      run's formatting onto each new run), matching the diff structure.

function Split-Run {
    param($TextRange, $OldText, $NewParts)

    $full = $TextRange.Text
    $idx = $full.IndexOf($OldText)
    if ($idx -lt 0) {
        throw "Could not find old text: $OldText"
    }

    $startPos = $idx + 1
    $combined = [string]::Join("", $NewParts)

    # Step 1: replace the whole old span with the full new text.
    $whole = $TextRange.Characters($startPos, $OldText.Length)
    $whole.Text = $combined

    # Step 2: re-set each part on its own sub-range so separate runs appear.
    $cursor = $startPos
    foreach ($part in $NewParts) {
        $sub = $TextRange.Characters($cursor, $part.Length)
        $sub.Text = $part
        $cursor += $part.Length
    }
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 4 - "Closet- Sneakerhead" (Dean)
# ---------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$tr4 = $s4.Shapes.Item(2).TextFrame.TextRange

Split-Run $tr4 "Dean has no way of knowing how much his sneakers are worth, does not have access to ppl willing to buy " @(
    "Dean has no way of knowing how much his sneakers are worth, does not have access to ",
    "people ",
    "willing to buy "
)

Split-Run $tr4 "Dean opens webpage goes to `"Sneakerhead`" interest group and enters into a convo with other sneakerheads to advice as to where to find the shoe" @(
    "Dean opens webpage goes to “Sneakerhead” interest group and enters into a convo with other sneakerheads to ",
    "get advice ",
    "as to where to find the shoe"
)

Split-Run $tr4 "Dean opens webpage and  creates a listing putting his show up for sell" @(
    "Dean opens webpage and  creates a listing putting ",
    "his, shoe up ",
    "for sell"
)

# ---------------------------------------------------------------------
# Slide 6 - "Closet- Purse Fashionista" (Shayla)
# ---------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$tr6 = $s6.Shapes.Item(2).TextFrame.TextRange

Split-Run $tr6 "Shayla is curious as to what is latest purse of Louis Vuitton" @(
    "Shayla is curious as to what ",
    "is the Louis Vuitton ",
    "latest purse "
)

Split-Run $tr6 "Shayla needs to find affordable access to other authenticate high end purses  " @(
    "Shayla needs to find affordable access to other ",
    "authentic ",
    "high end purses  "
)

# ---------------------------------------------------------------------
# Slide 8 - "Closet- Gadget Geek" (Marc)
# ---------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$tr8 = $s8.Shapes.Item(2).TextFrame.TextRange

Split-Run $tr8 "Character and market:" @(
    "Character and ",
    "market:"
)

Split-Run $tr8 "Marc owns a pair of Google Glass that he doesn't mind to rent out to friends for several days and make some money. " @(
    "Marc owns a pair of Google ",
    "Glasses ",
    "that he doesn’t mind to rent out to friends for several days and make some money. "
)

Split-Run $tr8 "Marc lists his Google Glass and can select the renter by checking his score, " @(
    "Marc lists his Google ",
    "Glasses ",
    "and can select the renter by checking his score, "
)

Split-Run $tr8 "facebook" @(
    "F",
    "acebook"
)
